$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = 390.0770000151745
$ws.Range("C5").Value = 80
$ws.Range("D5").Value = 38.7457683870282
$ws.Range("E5").Value = 131.5638295156349

$ws.Range("B6").Value = 390.0770000151745
$ws.Range("C6").Value = 120
$ws.Range("D6").Value = 58.1186525805423
$ws.Range("E6").Value = 197.3457442734524

$ws.Range("B7").Value = 391.0119755478833
$ws.Range("C7").Value = 39.99999999999793
$ws.Range("D7").Value = 19.37288419351511
$ws.Range("E7").Value = 65.62461911308367

$ws.Range("B8").Value = 391.0119755478833
$ws.Range("C8").Value = 59.99999999999687
$ws.Range("D8").Value = 29.05932629027266
$ws.Range("E8").Value = 98.4369286696255

$ws.Range("B9").Value = 390.3814865225988
$ws.Range("C9").Value = 39.99999999999871
$ws.Range("D9").Value = [double]"1.636578872421524E-12"
$ws.Range("E9").Value = 59.15754605399023

$ws.Range("B10").Value = 383.7093736109052
$ws.Range("C10").Value = 200
$ws.Range("D10").Value = 96.86442096757048
$ws.Range("E10").Value = 334.3678018927355

$ws.Range("B11").Value = 390.7147973406008
$ws.Range("C11").Value = 150
$ws.Range("D11").Value = 72.64831572567788
$ws.Range("E11").Value = 246.2794998292696

$ws.Range("B15").Value = 0
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0

$ws.Range("B16").Value = 1354.966545725536
$ws.Range("C16").Value = 20.01755005170917
$ws.Range("D16").Value = -0.0001865144429667751
$ws.Range("E16").Value = 14.77346370268729

$ws.Range("B17").Value = 391.0120092064996
$ws.Range("C17").Value = 49.9999999999984
$ws.Range("D17").Value = [double]"2.041815605480224E-12"
$ws.Range("E17").Value = 127.8733103401842

$ws.Range("B18").Value = 586.0722512117063
$ws.Range("C18").Value = 30
$ws.Range("D18").Value = [double]"7.958078640513122E-16"
$ws.Range("E18").Value = 51.18822796673089

$ws.Range("B22").Value = 1354.966545725536
$ws.Range("C22").Value = 81.52554998775024
$ws.Range("D22").Value = -0.0007596180594018734
$ws.Range("E22").Value = 60.16794076037835

$ws.Range("B23").Value = 586.5178881786698
$ws.Range("C23").Value = -9.99999999999968
$ws.Range("D23").Value = [double]"-4.078231086168671E-13"
$ws.Range("E23").Value = 17.0497783640551

$ws.Range("B24").Value = 586.0722512124808
$ws.Range("C24").Value = 50
$ws.Range("D24").Value = [double]"1.818989403545857E-15"
$ws.Range("E24").Value = 85.31371327777207

$ws.Range("B25").Value = 602.9735305280847
$ws.Range("C25").Value = 100.0000000000006
$ws.Range("D25").Value = [double]"1.33240973809734E-13"
$ws.Range("E25").Value = 165.8447592424503

$ws.Range("B29").Value = 677.4830753872644
$ws.Range("C29").Value = -122.2883249816272
$ws.Range("D29").Value = 0.001139429630412053
$ws.Range("E29").Value = 180.5038272445742

$ws.Range("B30").Value = 677.4830753872644
$ws.Range("C30").Value = -122.2883249816272
$ws.Range("D30").Value = 0.001139429630412053
$ws.Range("E30").Value = 180.5038272445742

$ws.Range("B31").Value = 586.5178881786698
$ws.Range("C31").Value = -29.99999999999904
$ws.Range("D31").Value = [double]"-1.224861989612691E-12"
$ws.Range("E31").Value = 51.1493350921653

$ws.Range("B32").Value = 586.0722512117063
$ws.Range("C32").Value = -49.99999999999999
$ws.Range("D32").Value = [double]"-8.526512829121202E-16"
$ws.Range("E32").Value = 85.31371327788482

$ws.Range("B33").Value = 602.9735305280847
$ws.Range("C33").Value = -40.00000000000024
$ws.Range("D33").Value = [double]"-5.312017492542508E-14"
$ws.Range("E33").Value = 66.33790369698012

$ws.Range("B37").Value = 2029.50666434999
$ws.Range("C37").Value = -140.0000000000007
$ws.Range("D37").Value = -90.43106644167099
$ws.Range("E37").Value = 47.41302236410993

$ws.Range("B38").Value = 2029.50666434999
$ws.Range("C38").Value = -140.0000000000007
$ws.Range("D38").Value = -79.34165632123961
$ws.Range("E38").Value = 45.77809055845096

$ws.Range("B42").Value = 586.0722512124808
$ws.Range("C42").Value = -10
$ws.Range("D42").Value = [double]"2.1316282072803E-16"
$ws.Range("E42").Value = 17.06274265555441

$ws.Range("B49").Value = 19991.83895448359
$ws.Range("C49").Value = 1993.104211821479
$ws.Range("D49").Value = 100.746558779756
$ws.Range("E49").Value = 100.6650747866874
$ws.Range("F49").Value = 85.7868672670006
$ws.Range("G49").Value = 55.46602179722338
$ws.Range("H49").Value = 3.821382138521454
$ws.Range("I49").Value = 33.29352170595222

$ws.Range("B50").Value = 19991.83895448359
$ws.Range("C50").Value = 1993.104211821479
$ws.Range("D50").Value = 100.746558779756
$ws.Range("E50").Value = 100.6650747866874
$ws.Range("F50").Value = 85.7868672670006
$ws.Range("G50").Value = 55.46602179722338
$ws.Range("H50").Value = 3.821382138521454
$ws.Range("I50").Value = 33.29352170595222

$ws.Range("B51").Value = 1993.104211821479
$ws.Range("C51").Value = 390.0770000151745
$ws.Range("D51").Value = 201.3301495654292
$ws.Range("E51").Value = 200.0000000000006
$ws.Range("F51").Value = 110.9320435944472
$ws.Range("G51").Value = 96.86442096757094
$ws.Range("H51").Value = 66.58704340988417
$ws.Range("I51").Value = 328.9095737890887

$ws.Range("B52").Value = 19992.41183508806
$ws.Range("C52").Value = 391.0119755478833
$ws.Range("D52").Value = 216.776119104398
$ws.Range("E52").Value = 215.2747887781741
$ws.Range("F52").Value = 114.0867981039403
$ws.Range("G52").Value = 98.96287412116075
$ws.Range("H52").Value = 7.074249419306123
$ws.Range("I52").Value = 349.8477544884278

$ws.Range("B53").Value = 19997.49870920183
$ws.Range("C53").Value = 392.6505533150195
$ws.Range("D53").Value = 311.293282072608
$ws.Range("E53").Value = 309.4268461673401
$ws.Range("F53").Value = 143.7101467674327
$ws.Range("G53").Value = 122.7980265354749
$ws.Range("H53").Value = 9.898946359505542
$ws.Range("I53").Value = 489.5050681257725

$ws.Range("B54").Value = 19997.49870920183
$ws.Range("C54").Value = 2009.911759163232
$ws.Range("D54").Value = -217.7307417301425
$ws.Range("E54").Value = -218.036026866227
$ws.Range("F54").Value = -135.0164491780118
$ws.Range("G54").Value = -167.2857969090567
$ws.Range("H54").Value = 7.39665112388319
$ws.Range("I54").Value = 78.94153801240478

$ws.Range("B58").Value = 1993.104211821479
$ws.Range("C58").Value = 0
$ws.Range("D58").Value = [double]"3.972464406649234E-09"
$ws.Range("E58").Value = 0
$ws.Range("F58").Value = [double]"1.058791184067875E-25"
$ws.Range("G58").Value = 0
$ws.Range("H58").Value = [double]"1.150719253218108E-09"
$ws.Range("I58").Value = 0

$ws.Range("B59").Value = 1993.104211821479
$ws.Range("C59").Value = 1354.966545725536
$ws.Range("D59").Value = [double]"3.972464406649234E-09"
$ws.Range("E59").Value = [double]"-1.835934499704411E-09"
$ws.Range("F59").Value = [double]"1.058791184067875E-25"
$ws.Range("G59").Value = [double]"-5.293955920339377E-26"
$ws.Range("H59").Value = [double]"1.150719253218108E-09"
$ws.Range("I59").Value = [double]"1.354966545725537E-09"

$ws.Range("B60").Value = 391.0119755478833
$ws.Range("C60").Value = 586.5178881786698
$ws.Range("D60").Value = 75.220734335877
$ws.Range("E60").Value = 75.19441776525738
$ws.Range("F60").Value = 50.46965844039665
$ws.Range("G60").Value = 50.44431608087888
$ws.Range("H60").Value = 133.7864627940851
$ws.Range("I60").Value = 154.3816362692046

$ws.Range("B61").Value = 390.7147973406008
$ws.Range("C61").Value = 586.0722512124808
$ws.Range("D61").Value = -25.16261565403326
$ws.Range("E61").Value = -25.16461410261487
$ws.Range("F61").Value = -50.44621212878806
$ws.Range("G61").Value = -50.44426516573817
$ws.Range("H61").Value = 83.28957374873137
$ws.Range("I61").Value = 96.17006688312448

$ws.Range("B62").Value = 390.7147973406008
$ws.Range("C62").Value = 586.0722512117063
$ws.Range("D62").Value = -19.99964045259352
$ws.Range("E62").Value = -20.0011205116319
$ws.Range("F62").Value = -0.004219132152220368
$ws.Range("G62").Value = -0.001926431546492438
$ws.Range("H62").Value = 29.54892875678961
$ws.Range("I62").Value = 34.12664753478998

$ws.Range("B63").Value = 2009.911759163232
$ws.Range("C63").Value = 602.9735305280847
$ws.Range("D63").Value = 59.99831629691234
$ws.Range("E63").Value = 59.9985584158129
$ws.Range("F63").Value = -0.004668075256308243
$ws.Range("G63").Value = 0.0009960316073326252
$ws.Range("H63").Value = 17.23487101640739
$ws.Range("I63").Value = 99.50396450329896

$ws.Range("B67").Value = 1354.966545725536
$ws.Range("C67").Value = 677.4830753872644
$ws.Range("D67").Value = -244.7914401695775
$ws.Range("E67").Value = -244.791451197598
$ws.Range("F67").Value = 0.000800023647319108
$ws.Range("G67").Value = 0.001738821151774914
$ws.Range("H67").Value = 180.6622819899379
$ws.Range("I67").Value = 361.3242919723191

$ws.Range("B68").Value = 586.5178881786698
$ws.Range("C68").Value = 391.0120092064996
$ws.Range("D68").Value = 50.00088414298119
$ws.Range("E68").Value = 50.00032187178126
$ws.Range("F68").Value = -0.003077345396932287
$ws.Range("G68").Value = 0.001435616236097076
$ws.Range("H68").Value = 85.25069180540693
$ws.Range("I68").Value = 127.8807575797848

$ws.Range("B72").Value = 19998.72825263309
$ws.Range("C72").Value = 19991.83895448359
$ws.Range("D72").Value = 201.5599338450356
$ws.Range("E72").Value = 201.4931175594023
$ws.Range("F72").Value = 62.34267563938607
$ws.Range("G72").Value = 171.5737345339758
$ws.Range("H72").Value = 6.090886458193824
$ws.Range("I72").Value = 7.642764277038633

$ws.Range("B73").Value = 19998.72825263309
$ws.Range("C73").Value = 19992.41183508806
$ws.Range("D73").Value = 216.838315631131
$ws.Range("E73").Value = 216.7761191043386
$ws.Range("F73").Value = 4.849509963305447
$ws.Range("G73").Value = 114.0867981039294
$ws.Range("H73").Value = 6.261562773531336
$ws.Range("I73").Value = 7.074249419303102

$ws.Range("B74").Value = 391.0119755478833
$ws.Range("C74").Value = 390.3814865225988
$ws.Range("D74").Value = 40.06456786865245
$ws.Range("E74").Value = 39.99999999999748
$ws.Range("F74").Value = 0.05249420216037322
$ws.Range("G74").Value = [double]"1.8712853488978E-12"
$ws.Range("H74").Value = 59.15754605399238
$ws.Range("I74").Value = 59.15754605399238

$ws.Range("B75").Value = 19998.72825263309
$ws.Range("C75").Value = 19997.49870920183
$ws.Range("D75").Value = 93.57035151509126
$ws.Range("E75").Value = 93.56254034252444
$ws.Range("F75").Value = -64.17226655444813
$ws.Range("G75").Value = 8.69369758948239
$ws.Range("H75").Value = 3.275554842425507
$ws.Range("I75").Value = 2.712936233227784

$ws.Range("B76").Value = 392.6505533150195
$ws.Range("C76").Value = 383.7093736109052
$ws.Range("D76").Value = 204.1254871459476
$ws.Range("E76").Value = 199.9999999999996
$ws.Range("F76").Value = 100.2184755577707
$ws.Range("G76").Value = 96.86442096756912
$ws.Range("H76").Value = 334.3678018927345
$ws.Range("I76").Value = 334.3678018927345

$ws.Range("B77").Value = 392.6505533150195
$ws.Range("C77").Value = 390.7147973406008
$ws.Range("D77").Value = 105.3013590213881
$ws.Range("E77").Value = 104.8387077912311
$ws.Range("F77").Value = 22.57955097770001
$ws.Range("G77").Value = 22.20341176619027
$ws.Range("H77").Value = 158.3617801789023
$ws.Range("I77").Value = 158.3617801789023

$ws.Range("B78").Value = 2009.911759163232
$ws.Range("C78").Value = 2029.50666434999
$ws.Range("D78").Value = -278.0360268661007
$ws.Range("E78").Value = -280.0000000000017
$ws.Range("F78").Value = -167.2791613907123
$ws.Range("G78").Value = -169.7727227629073
$ws.Range("H78").Value = 93.20697572347362
$ws.Range("I78").Value = 93.15205006387311

$ws.Range("B82").Value = 586.0722512124808
$ws.Range("C82").Value = 586.5178881786698
$ws.Range("D82").Value = -65.16158679284047
$ws.Range("E82").Value = -65.21022090988882
$ws.Range("F82").Value = -50.44377331627581
$ws.Range("G82").Value = -50.4833059055368
$ws.Range("H82").Value = 140.6007605085263
$ws.Range("I82").Value = 140.6007605085263
